$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint1")

# Row 2
$ws.Range("E2").Value = 25
$ws.Range("F2").Value = 3

# Row 3
$ws.Range("E3").Value = 25
$ws.Range("F3").Value = 3

# Row 4 - status changes to COMPLETED
$ws.Range("D4").Value = "COMPLETED"
$ws.Range("E4").Value = 20
$ws.Range("F4").Value = 2

# Row 5
$ws.Range("E5").Value = 22
$ws.Range("F5").Value = 2

# Row 6 - status changes to COMPLETED
$ws.Range("D6").Value = "COMPLETED"
$ws.Range("E6").Value = 20
$ws.Range("F6").Value = 3

# Row 7
$ws.Range("E7").Value = 20
$ws.Range("F7").Value = 2

# Row 8
$ws.Range("E8").Value = 15
$ws.Range("F8").Value = 2

# Row 9
$ws.Range("E9").Value = 20
$ws.Range("F9").Value = 2

# Row 10
$ws.Range("E10").Value = 20
$ws.Range("F10").Value = 2

# Row 11
$ws.Range("E11").Value = 20
$ws.Range("F11").Value = 2

# Clear the stale "coding right now" / "Jim rowland" notes
$ws.Range("B21").Value = ""
$ws.Range("B22").Value = ""

$ws.Range("D2").Select()
